# Swap the contents of columns C (codeforiati:group-code) and D
# (codeforiati:group-name) for every row of data, including the header row.
#
# The underlying OOXML diff shows that, for each row, the shared-string
# entries used by columns C and D were exchanged - i.e. the two columns'
# values simply swap places throughout the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value()
    $dValue = $dCell.Value()

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
